$d = $word.ActiveDocument

# 1) Consolidate the "My C++ skills..." paragraph into a single run and
#    append the new "iterating on gameplay systems..." sentence fragment.
$p13 = $d.Paragraphs.Item(13)
$r13 = $p13.Range
[void]$r13.MoveEnd(1, -1)   # exclude the paragraph mark from the range
$newText = "My C++ skills are at their peak through rigorous practice with the use of pointers and a better understanding of Data Structures from the Collision System and Memory Manager that I created. I love delving into 3D Math and am relearning it in a better way, with a heavy focus on understanding it through geometry and visualizing it, for use specifically in games. The Action games that I have worked on and am currently working on have given me experience in bringing the design, engineering, art and animation in them together, and in collaborating with and learning from the people involved in them, as well as iterating on gameplay systems to get them to their best possible form for the game."
$r13.Text = $newText

# 2) Move the (hidden) "_GoBack" bookmark from the end of the third body
#    paragraph to the very start of the "I love Action Games..." paragraph.
$p14 = $d.Paragraphs.Item(14)
$pos = $p14.Range.Start
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
